$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Update summary header values ---
# Valor Mora total
$ws.Range("E11").Value = 487485
# Cant. Trabajadores
$ws.Range("C13").Value = 5
# Cant. Periodos
$ws.Range("F13").Value = 8

# --- Remove the now-obsolete data rows ---
# Keep rows 16-23 (normal style) untouched, and delete rows 24-31 so that
# the old row 32 (which carries the special "last row" bottom-border
# style) shifts up to become the new row 24. This also shifts the footer
# rows (old 37-38) up to 29-30 and updates dimension / merged cell
# references automatically.
$ws.Range("A24:A31").EntireRow.Delete()

# --- Overwrite the remaining data rows (16-24) with the new dataset ---
$ws.Range("B16").Value = "CE"
$ws.Range("C16").Value = "668215"
$ws.Range("D16").Value = "JOSE RAMON LEZAMA DIAZ"
$ws.Range("E16").Value = "1703"
$ws.Range("F16").Value = 7869
$ws.Range("G16").Value = 737717

$ws.Range("B17").Value = "CC"
$ws.Range("C17").Value = "9294722"
$ws.Range("D17").Value = "VICTOR MANUEL BENITEZ MONTIEL"
$ws.Range("E17").Value = "1704"
$ws.Range("F17").Value = 29509
$ws.Range("G17").Value = 1700000

$ws.Range("B18").Value = "CE"
$ws.Range("C18").Value = "668215"
$ws.Range("D18").Value = "JOSE RAMON LEZAMA DIAZ"
$ws.Range("E18").Value = "1704"
$ws.Range("F18").Value = 29509
$ws.Range("G18").Value = 737717

$ws.Range("B19").Value = "CC"
$ws.Range("C19").Value = "1127585376"
$ws.Range("D19").Value = "MILEIS DE JESUS CONEO ALVAREZ"
$ws.Range("E19").Value = "1705"
$ws.Range("F19").Value = 15738
$ws.Range("G19").Value = 737717

$ws.Range("B20").Value = "CC"
$ws.Range("C20").Value = "9294722"
$ws.Range("D20").Value = "VICTOR MANUEL BENITEZ MONTIEL"
$ws.Range("E20").Value = "1707"
$ws.Range("F20").Value = 29509
$ws.Range("G20").Value = 1700000

$ws.Range("B21").Value = "CC"
$ws.Range("C21").Value = "9294722"
$ws.Range("D21").Value = "VICTOR MANUEL BENITEZ MONTIEL"
$ws.Range("E21").Value = "1708"
$ws.Range("F21").Value = 29509
$ws.Range("G21").Value = 1700000

$ws.Range("B22").Value = "CC"
$ws.Range("C22").Value = "9294722"
$ws.Range("D22").Value = "VICTOR MANUEL BENITEZ MONTIEL"
$ws.Range("E22").Value = "1709"
$ws.Range("F22").Value = 29509
$ws.Range("G22").Value = 1700000

$ws.Range("B23").Value = "CC"
$ws.Range("C23").Value = "1023165078"
$ws.Range("D23").Value = "CARMEN EVITA ANGEL MARTON"
$ws.Range("E23").Value = "2003"
$ws.Range("F23").Value = 315000
$ws.Range("G23").Value = 11250000

$ws.Range("B24").Value = "CC"
$ws.Range("C24").Value = "1143359392"
$ws.Range("D24").Value = "RUBEN DARIO CARAZO SEQUEA"
$ws.Range("E24").Value = "2111"
$ws.Range("F24").Value = 1333
$ws.Range("G24").Value = 1000000
